$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1603853333333333
$ws.Range("H2").Value = 0.481156
$ws.Range("I2").Value = 0.01032935781992836
$ws.Range("J2").Value = 0.01042870175281933
$ws.Range("M2").Value = 1.037532
$ws.Range("N2").Value = 3.112596
$ws.Range("O2").Value = 0.04166450179684251
$ws.Range("P2").Value = 0.0439159257402554
$ws.Range("Q2").Value = 0.166404915664
$ws.Range("R2").Value = 1.497644240976
$ws.Range("S2").Value = 0.0004303675474486344
$ws.Range("T2").Value = 0.0004579860917440852

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1603853333333333
$ws.Range("H3").Value = 0.481156
$ws.Range("I3").Value = 0.01032935781992836
$ws.Range("J3").Value = 0.01042870175281933
$ws.Range("O3").Value = 0.2561129158441639
$ws.Range("P3").Value = 0.2699524849277078
$ws.Range("Q3").Value = 1.022895902351555
$ws.Range("R3").Value = 9.206063121163998
$ws.Range("S3").Value = 0.002645481950059569
$ws.Range("T3").Value = 0.002815253952743521

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1603853333333333
$ws.Range("H4").Value = 0.481156
$ws.Range("I4").Value = 0.01032935781992836
$ws.Range("J4").Value = 0.01042870175281933
$ws.Range("M4").Value = 6.239319333333333
$ws.Range("N4").Value = 18.717958
$ws.Range("O4").Value = 0.2505543265891952
$ws.Range("P4").Value = 0.2640935262839185
$ws.Range("Q4").Value = 1.000695311049778
$ws.Range("R4").Value = 9.006257799447999
$ws.Range("S4").Value = 0.002588065292670988
$ws.Range("T4").Value = 0.002754152620465339

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1603853333333333
$ws.Range("H5").Value = 0.481156
$ws.Range("I5").Value = 0.01032935781992836
$ws.Range("J5").Value = 0.01042870175281933
$ws.Range("M5").Value = 3.8299385
$ws.Range("N5").Value = 7.659877
$ws.Range("O5").Value = 0.1538000558200097
$ws.Range("P5").Value = 0.1080739644693659
$ws.Range("Q5").Value = 0.6142659629686666
$ws.Range("R5").Value = 3.685595777812
$ws.Range("S5").Value = 0.001588655809289835
$ws.Range("T5").Value = 0.00112707114269581

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1603853333333333
$ws.Range("H6").Value = 0.481156
$ws.Range("I6").Value = 0.01032935781992836
$ws.Range("J6").Value = 0.01042870175281933
$ws.Range("M6").Value = 7.417532333333334
$ws.Range("N6").Value = 22.252597
$ws.Range("O6").Value = 0.2978681999497886
$ws.Range("P6").Value = 0.3139640985787523
$ws.Range("Q6").Value = 1.189663395792445
$ws.Range("R6").Value = 10.706970562132
$ws.Range("S6").Value = 0.003076787220459334
$ws.Range("T6").Value = 0.003274237945170577

$ws.Range("I7").Value = 0.9610926076617912
$ws.Range("J7").Value = 0.970336039943066
$ws.Range("M7").Value = 1.037532
$ws.Range("N7").Value = 3.112596
$ws.Range("O7").Value = 0.04166450179684251
$ws.Range("P7").Value = 0.0439159257402554
$ws.Range("Q7").Value = 15.483105253136
$ws.Range("R7").Value = 139.347947278224
$ws.Range("S7").Value = 0.04004344467885675
$ws.Range("T7").Value = 0.04261320547323319

$ws.Range("I8").Value = 0.9610926076617912
$ws.Range("J8").Value = 0.970336039943066
$ws.Range("O8").Value = 0.2561129158441639
$ws.Range("P8").Value = 0.2699524849277078
$ws.Range("S8").Value = 0.2461482301445324
$ws.Range("T8").Value = 0.2619446251975422

$ws.Range("I9").Value = 0.9610926076617912
$ws.Range("J9").Value = 0.970336039943066
$ws.Range("M9").Value = 6.239319333333333
$ws.Range("N9").Value = 18.717958
$ws.Range("O9").Value = 0.2505543265891952
$ws.Range("P9").Value = 0.2640935262839185
$ws.Range("Q9").Value = 93.10945392135021
$ws.Range("R9").Value = 837.9850852921519
$ws.Range("S9").Value = 0.2408059111025537
$ws.Range("T9").Value = 0.2562594664689375

$ws.Range("I10").Value = 0.9610926076617912
$ws.Range("J10").Value = 0.970336039943066
$ws.Range("M10").Value = 3.8299385
$ws.Range("N10").Value = 7.659877
$ws.Range("O10").Value = 0.1538000558200097
$ws.Range("P10").Value = 0.1080739644693659
$ws.Range("Q10").Value = 57.15422840793133
$ws.Range("R10").Value = 342.925370447588
$ws.Range("S10").Value = 0.1478160967065822
$ws.Range("T10").Value = 0.1048680627041521

$ws.Range("I11").Value = 0.9610926076617912
$ws.Range("J11").Value = 0.970336039943066
$ws.Range("M11").Value = 7.417532333333334
$ws.Range("N11").Value = 22.252597
$ws.Range("O11").Value = 0.2978681999497886
$ws.Range("P11").Value = 0.3139640985787523
$ws.Range("Q11").Value = 110.6919438008076
$ws.Range("R11").Value = 996.2274942072679
$ws.Range("S11").Value = 0.2862789250292662
$ws.Range("T11").Value = 0.304650680099201

$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.5
$ws.Range("G12").Value = 0.443735
$ws.Range("H12").Value = 0.88747
$ws.Range("I12").Value = 0.02857803451828042
$ws.Range("J12").Value = 0.01923525830411462
$ws.Range("M12").Value = 1.037532
$ws.Range("N12").Value = 3.112596
$ws.Range("O12").Value = 0.04166450179684251
$ws.Range("P12").Value = 0.0439159257402554
$ws.Range("Q12").Value = 0.46038926202
$ws.Range("R12").Value = 2.76233557212
$ws.Range("S12").Value = 0.001190689570537122
$ws.Range("T12").Value = 0.0008447341752781289

$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.5
$ws.Range("G13").Value = 0.443735
$ws.Range("H13").Value = 0.88747
$ws.Range("I13").Value = 0.02857803451828042
$ws.Range("J13").Value = 0.01923525830411462
$ws.Range("O13").Value = 0.2561129158441639
$ws.Range("P13").Value = 0.2699524849277078
$ws.Range("Q13").Value = 2.830026310988333
$ws.Range("R13").Value = 16.98015786593
$ws.Range("S13").Value = 0.007319203749571963
$ws.Range("T13").Value = 0.005192605777422069

$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.5
$ws.Range("G14").Value = 0.443735
$ws.Range("H14").Value = 0.88747
$ws.Range("I14").Value = 0.02857803451828042
$ws.Range("J14").Value = 0.01923525830411462
$ws.Range("M14").Value = 6.239319333333333
$ws.Range("N14").Value = 18.717958
$ws.Range("O14").Value = 0.2505543265891952
$ws.Range("P14").Value = 0.2640935262839185
$ws.Range("Q14").Value = 2.768604364376666
$ws.Range("R14").Value = 16.61162618626
$ws.Range("S14").Value = 0.007160350193970525
$ws.Range("T14").Value = 0.005079907194515656

$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.5
$ws.Range("G15").Value = 0.443735
$ws.Range("H15").Value = 0.88747
$ws.Range("I15").Value = 0.02857803451828042
$ws.Range("J15").Value = 0.01923525830411462
$ws.Range("M15").Value = 3.8299385
$ws.Range("N15").Value = 7.659877
$ws.Range("O15").Value = 0.1538000558200097
$ws.Range("P15").Value = 0.1080739644693659
$ws.Range("Q15").Value = 1.6994777602975
$ws.Range("R15").Value = 6.79791104119
$ws.Range("S15").Value = 0.004395303304137692
$ws.Range("T15").Value = 0.002078830622517958

$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.5
$ws.Range("G16").Value = 0.443735
$ws.Range("H16").Value = 0.88747
$ws.Range("I16").Value = 0.02857803451828042
$ws.Range("J16").Value = 0.01923525830411462
$ws.Range("M16").Value = 7.417532333333334
$ws.Range("N16").Value = 22.252597
$ws.Range("O16").Value = 0.2978681999497886
$ws.Range("P16").Value = 0.3139640985787523
$ws.Range("Q16").Value = 3.291418709931667
$ws.Range("R16").Value = 19.74851225959
$ws.Range("S16").Value = 0.008512487700063113
$ws.Range("T16").Value = 0.006039180534380809
